$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: clear all columns except A, B, keep C (renamed) ---
$ws.Range("C1").Value = "temp"
$ws.Range("D1:R1").ClearContents() | Out-Null

# --- Row 2: update C2, clear the rest of the extra/date columns ---
$ws.Range("C2").Value = 3
$ws.Range("D2:R2").ClearContents() | Out-Null

# --- New rows 3 and 4 ---
$ws.Range("A3").Value = "asdf"
$ws.Range("C3").Value = 4
$ws.Range("A4").Value = "qwer"
$ws.Range("C4").Value = 5

# --- Hyperlinks: rebuild B2, B3, B4 so that the final order/ids match ---
$ws.Range("B2").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:wonjae@example.com", [Type]::Missing, [Type]::Missing, "mailto:wonjae@example.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:wonjae@example.com", [Type]::Missing, [Type]::Missing, "mailto:wonjae@example.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:wonjae@example.com", [Type]::Missing, [Type]::Missing, "mailto:wonjae@example.com") | Out-Null

$ws.Range("B4").Value = "qwer"
$ws.Range("B3").Value = "asdf"
$ws.Range("B2").Value = "wonjae@example.com"

$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# --- Selection ---
$ws.Range("I22").Select() | Out-Null
